# A new weekly price record was inserted into the data table.
# It is inserted as a new row right after the existing row 513, which
# shifts every following row (514..541) down by one (515..542), and the
# new row 514 receives an exact copy of what row 513 used to contain.
# Row 513 itself is then updated with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 514, pushing old rows 514-541 down to 515-542.
$ws.Rows.Item(514).Insert()

# The new row 514 takes on the data that row 513 used to hold
# (only the used columns A:T, to avoid touching the rest of the row).
$ws.Range("A513:T513").Copy($ws.Range("A514:T514"))

# Row 513 now becomes the new price record.
$ws.Range("D513").Value = 44753
$ws.Range("M513").Value = 400
$ws.Range("N513").Value = 26000
$ws.Range("O513").Value = 27000
$ws.Range("P513").Value = 26500
$ws.Range("S513").Value = 1325
